# Decrement the "剩余" (remaining) value in column E by 1 for every data
# row (rows 2-99), except row 36 which is left unchanged, matching the
# author's commit "自动更新Excel文件 - Sat Nov 15 23:23:16 UTC 2025".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)   # column E
    $current = $cell.Value2
    if ($null -ne $current) {
        $cell.Value2 = $current - 1
    }
}
